# compute the fraction BB:BP
# The "Data" sheet maps Eora_classification -> Sector -> Industry_category.
# Rows 17 and 18 were mis-tagged with Sector "Business" (a label that isn't
# used anywhere else in the sheet); both should instead be tagged "Services",
# matching the rest of that Industry_category group. Once "Business" is no
# longer referenced by any cell, Excel drops it from the shared-strings table
# on save, which is why all later shared-string indices shift down by one -
# that's an automatic side effect of the text fix below, not something to do
# by hand.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")
$ws.Activate()

$ws.Range("B17").Value = "Services"
$ws.Range("B18").Value = "Services"

# Scroll the sheet so row 9 is at the top (best-effort view-state update).
[void]$ws.Range("A9").Select()
$excel.ActiveWindow.ScrollRow = 9
$excel.ActiveWindow.ScrollColumn = 1

# Restore the original selection/active cell on the sheet.
[void]$ws.Range("B15").Select()
